$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.665.32'
$ws.Range("E2").Value = '  +3.53%  '

$ws.Range("D3").Value = '2.252.95'
$ws.Range("E3").Value = '  +3.49%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.635'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.89'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.05%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.646'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.52'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0964'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.24%  '

$ws.Range("E14").Value = '  +0.63%  '

$ws.Range("D15").Value = '2.589.50'
$ws.Range("E15").Value = '  +3.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.886'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.76%  '

$ws.Range("D18").Value = '2.250.77'
$ws.Range("E18").Value = '  +2.64%  '

$ws.Range("D19").Value = '42.683.98'
$ws.Range("E19").Value = '  +3.77%  '

$ws.Range("E20").Value = '  +3.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.22%  '

$ws.Range("E25").Value = '  +1.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.69%  '

$ws.Range("E30").Value = '  +1.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("E32").Value = '  +2.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.128'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +14.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0786'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.37%  '

$ws.Range("E36").Value = '  +1.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.64%  '

$ws.Range("E40").Value = '  +5.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.29%  '

$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.67%  '

$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.81%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.73%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.203'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.17%  '

$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.102'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.59%  '

$ws.Range("E49").Value = '  +5.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.49%  '
